# Update the "想去人数" (want-to-go count) values from 8 to 9
# for the two events on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9
    $ws.Range("F3").Value = 9
}
